$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: update existing scalar cell values across year-range blocks ---
$ws.Range("G2").Value = 1.06485605843523
$ws.Range("G3").Value = 0.0224147339421163
$ws.Range("G4").Value = 0.0224147339421163
$ws.Range("F9").Value = 0.04708
$ws.Range("G9").Value = 0.0824553401260571
$ws.Range("F10").Value = 0.04708
$ws.Range("G10").Value = 0.0824553401260571
$ws.Range("G11").Value = 0.573454166666667
$ws.Range("G12").Value = 0.573454166666667
$ws.Range("G15").Value = 1.16644360100518
$ws.Range("G16").Value = 0.0212635372853411
$ws.Range("G17").Value = 0.0212635372853411
$ws.Range("G22").Value = 0.0953865327301012
$ws.Range("G23").Value = 0.0953865327301012
$ws.Range("G26").Value = 0.5748779661016949
$ws.Range("G27").Value = 0.5748779661016949
$ws.Range("G33").Value = 0.0189584525395784
$ws.Range("G34").Value = 0.0189584525395784
$ws.Range("F43").Value = 0.5347
$ws.Range("F44").Value = 0.5347
$ws.Range("G49").Value = 1.35981589421305
$ws.Range("G50").Value = 0.0184684995019244
$ws.Range("G51").Value = 0.0184684995019244
$ws.Range("F60").Value = 0.5244
$ws.Range("G60").Value = 0.528237288135593
$ws.Range("I60").Value = 1.06891
$ws.Range("L60").Value = 0.29785
$ws.Range("N60").Value = 0.98625
$ws.Range("F61").Value = 0.5244
$ws.Range("G61").Value = 0.528237288135593
$ws.Range("I61").Value = 1.06891
$ws.Range("L61").Value = 0.29785
$ws.Range("N61").Value = 0.98625
$ws.Range("G66").Value = 1.49458443726334
$ws.Range("G67").Value = 0.0160617198409075
$ws.Range("G68").Value = 0.0160617198409075
$ws.Range("F77").Value = 0.5122
$ws.Range("G77").Value = 0.528041666666667
$ws.Range("I77").Value = 1.06435
$ws.Range("L77").Value = 0.29785
$ws.Range("M77").Value = 0.81453
$ws.Range("N77").Value = 0.98598
$ws.Range("F78").Value = 0.5122
$ws.Range("G78").Value = 0.528041666666667
$ws.Range("I78").Value = 1.06435
$ws.Range("L78").Value = 0.29785
$ws.Range("M78").Value = 0.81453
$ws.Range("N78").Value = 0.98598
$ws.Range("G83").Value = 1.42032774979583
$ws.Range("G84").Value = 0.0154084521615795
$ws.Range("G85").Value = 0.0154084521615795
$ws.Range("G86").Value = 831.0793090734099
$ws.Range("G87").Value = 831.0793090734099
$ws.Range("G88").Value = 831.0793090734099
$ws.Range("G89").Value = 831.0793090734099
$ws.Range("G90").Value = 0.123069706909843
$ws.Range("G91").Value = 0.123069706909843
$ws.Range("F94").Value = 0.49955
$ws.Range("G94").Value = 0.533175
$ws.Range("L94").Value = 0.29235
$ws.Range("M94").Value = 0.81453
$ws.Range("N94").Value = 1.00817
$ws.Range("F95").Value = 0.49955
$ws.Range("G95").Value = 0.533175
$ws.Range("L95").Value = 0.29235
$ws.Range("M95").Value = 0.81453
$ws.Range("N95").Value = 1.00817
$ws.Range("G100").Value = 1.18807144691514
$ws.Range("G101").Value = 0.0170416475749093
$ws.Range("G102").Value = 0.0170416475749093
$ws.Range("G103").Value = 808.233195667874
$ws.Range("G104").Value = 808.233195667874
$ws.Range("G105").Value = 808.233195667874
$ws.Range("G106").Value = 808.233195667874
$ws.Range("G107").Value = 0.104487393092963
$ws.Range("G108").Value = 0.104487393092963
$ws.Range("F111").Value = 0.4917
$ws.Range("G111").Value = 0.549875
$ws.Range("L111").Value = 0.2798
$ws.Range("M111").Value = 0.89272
$ws.Range("N111").Value = 1.01609
$ws.Range("F112").Value = 0.4917
$ws.Range("G112").Value = 0.549875
$ws.Range("L112").Value = 0.2798
$ws.Range("M112").Value = 0.89272
$ws.Range("N112").Value = 1.01609
$ws.Range("G117").Value = 0.940235467399761
$ws.Range("G118").Value = 0.0186992751937709
$ws.Range("G119").Value = 0.0186992751937709
$ws.Range("G120").Value = 1719.75861939669
$ws.Range("G121").Value = 1719.75861939669
$ws.Range("G122").Value = 1719.75861939669
$ws.Range("G123").Value = 1719.75861939669
$ws.Range("G124").Value = 0.102548212332608
$ws.Range("G125").Value = 0.102548212332608
$ws.Range("F128").Value = 0.5553
$ws.Range("G128").Value = 0.578225
$ws.Range("I128").Value = 1.07435
$ws.Range("L128").Value = 0.29235
$ws.Range("M128").Value = 0.89272
$ws.Range("N128").Value = 1.00208
$ws.Range("F129").Value = 0.5553
$ws.Range("G129").Value = 0.578225
$ws.Range("I129").Value = 1.07435
$ws.Range("L129").Value = 0.29235
$ws.Range("M129").Value = 0.89272
$ws.Range("N129").Value = 1.00208
$ws.Range("G135").Value = 0.0192755463802115
$ws.Range("G136").Value = 0.0192755463802115
$ws.Range("G137").Value = 1747.36878888821
$ws.Range("G138").Value = 1747.36878888821
$ws.Range("G139").Value = 1747.36878888821
$ws.Range("G140").Value = 1747.36878888821
$ws.Range("G141").Value = 0.125170744585316
$ws.Range("G142").Value = 0.125170744585316
$ws.Range("G145").Value = 0.602725
$ws.Range("M145").Value = 0.89272
$ws.Range("G146").Value = 0.602725
$ws.Range("M146").Value = 0.89272
$ws.Range("G151").Value = 0.768092306371656
$ws.Range("G152").Value = 0.0192585972276692
$ws.Range("G153").Value = 0.0192585972276692
$ws.Range("G154").Value = 1779.23721628284
$ws.Range("G155").Value = 1779.23721628284
$ws.Range("G156").Value = 1779.23721628284
$ws.Range("G157").Value = 1779.23721628284
$ws.Range("G158").Value = 0.110651997877846
$ws.Range("G159").Value = 0.110651997877846
$ws.Range("G168").Value = 0.699895440176964
$ws.Range("G169").Value = 0.0196145459911424
$ws.Range("G170").Value = 0.0196145459911424
$ws.Range("G175").Value = 0.110964078321441
$ws.Range("G176").Value = 0.110964078321441
$ws.Range("G185").Value = 0.690537832031756
$ws.Range("G186").Value = 0.018850909627506
$ws.Range("G187").Value = 0.018850909627506
$ws.Range("G188").Value = 3584.28693503237
$ws.Range("I188").Value = 16943.44808
$ws.Range("N188").Value = 12096.6776
$ws.Range("G189").Value = 3584.28693503237
$ws.Range("I189").Value = 16943.44808
$ws.Range("N189").Value = 12096.6776
$ws.Range("G190").Value = 3584.28693503237
$ws.Range("I190").Value = 16943.44808
$ws.Range("N190").Value = 12096.6776
$ws.Range("G191").Value = 3584.28693503237
$ws.Range("I191").Value = 16943.44808
$ws.Range("N191").Value = 12096.6776
$ws.Range("G192").Value = 0.116608580448276
$ws.Range("G193").Value = 0.116608580448276

# --- Part 2: append new "2019 - 2023" block of 17 rows (rows 202-218) ---
# Row 202
$ws.Cells.Item(202, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(202, 2).Value = 'Visual Clarity (Sediment class 4)'
$ws.Cells.Item(202, 3).Value = 'D'
$ws.Cells.Item(202, 4).Value = '2019 - 2023'
$ws.Cells.Item(202, 5).Value = 'RepSite'
$ws.Cells.Item(202, 6).Value = 0.325
$ws.Cells.Item(202, 7).Value = 0.670507690642522
$ws.Cells.Item(202, 8).Value = 3.77
$ws.Cells.Item(202, 9).Value = 2.42
$ws.Cells.Item(202, 12).Value = 0.36
$ws.Cells.Item(202, 13).Value = 1.5404
$ws.Cells.Item(202, 14).Value = 1.9856
$ws.Cells.Item(202, 15).Value = 1816987.417
$ws.Cells.Item(202, 16).Value = 5524893.755
$ws.Cells.Item(202, 17).Value = 'Palmerston North City'
$ws.Cells.Item(202, 18).Value = 'Manawatū'
$ws.Cells.Item(202, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(202, 20).Value = 'Mana_11a'
$ws.Cells.Item(202, 21).Value = 'm'
# Row 203
$ws.Cells.Item(203, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(203, 2).Value = 'DRP (95th Percentile)'
$ws.Cells.Item(203, 3).Value = 'C'
$ws.Cells.Item(203, 4).Value = '2019 - 2023'
$ws.Cells.Item(203, 5).Value = 'RepSite'
$ws.Cells.Item(203, 6).Value = 0.019
$ws.Cells.Item(203, 7).Value = 0.0191818181818182
$ws.Cells.Item(203, 8).Value = 0.04
$ws.Cells.Item(203, 9).Value = 0.0335
$ws.Cells.Item(203, 12).Value = 0.016
$ws.Cells.Item(203, 13).Value = 0.02615
$ws.Cells.Item(203, 14).Value = 0.0293
$ws.Cells.Item(203, 15).Value = 1816987.417
$ws.Cells.Item(203, 16).Value = 5524893.755
$ws.Cells.Item(203, 17).Value = 'Palmerston North City'
$ws.Cells.Item(203, 18).Value = 'Manawatū'
$ws.Cells.Item(203, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(203, 20).Value = 'Mana_11a'
$ws.Cells.Item(203, 21).Value = 'mg/L'
# Row 204
$ws.Cells.Item(204, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(204, 2).Value = 'DRP (Median)'
$ws.Cells.Item(204, 3).Value = 'D'
$ws.Cells.Item(204, 4).Value = '2019 - 2023'
$ws.Cells.Item(204, 5).Value = 'RepSite'
$ws.Cells.Item(204, 6).Value = 0.019
$ws.Cells.Item(204, 7).Value = 0.0191818181818182
$ws.Cells.Item(204, 8).Value = 0.04
$ws.Cells.Item(204, 9).Value = 0.0335
$ws.Cells.Item(204, 12).Value = 0.016
$ws.Cells.Item(204, 13).Value = 0.02615
$ws.Cells.Item(204, 14).Value = 0.0293
$ws.Cells.Item(204, 15).Value = 1816987.417
$ws.Cells.Item(204, 16).Value = 5524893.755
$ws.Cells.Item(204, 17).Value = 'Palmerston North City'
$ws.Cells.Item(204, 18).Value = 'Manawatū'
$ws.Cells.Item(204, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(204, 20).Value = 'Mana_11a'
$ws.Cells.Item(204, 21).Value = 'mg/L'
# Row 205
$ws.Cells.Item(205, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(205, 2).Value = 'E coli (>260)'
$ws.Cells.Item(205, 3).Value = 'E'
$ws.Cells.Item(205, 4).Value = '2019 - 2023'
$ws.Cells.Item(205, 5).Value = 'RepSite'
$ws.Cells.Item(205, 6).Value = 490
$ws.Cells.Item(205, 7).Value = 2624.66395229482
$ws.Cells.Item(205, 8).Value = 46000
$ws.Cells.Item(205, 9).Value = 12475.08199
$ws.Cells.Item(205, 10).Value = 47.2727272727273
$ws.Cells.Item(205, 11).Value = 60
$ws.Cells.Item(205, 12).Value = 380
$ws.Cells.Item(205, 13).Value = 2760
$ws.Cells.Item(205, 14).Value = 8250.5
$ws.Cells.Item(205, 15).Value = 1816987.417
$ws.Cells.Item(205, 16).Value = 5524893.755
$ws.Cells.Item(205, 17).Value = 'Palmerston North City'
$ws.Cells.Item(205, 18).Value = 'Manawatū'
$ws.Cells.Item(205, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(205, 20).Value = 'Mana_11a'
$ws.Cells.Item(205, 21).Value = '% exceedances over 260/100 mL'
# Row 206
$ws.Cells.Item(206, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(206, 2).Value = 'E coli (>540)'
$ws.Cells.Item(206, 3).Value = 'E'
$ws.Cells.Item(206, 4).Value = '2019 - 2023'
$ws.Cells.Item(206, 5).Value = 'RepSite'
$ws.Cells.Item(206, 6).Value = 490
$ws.Cells.Item(206, 7).Value = 2624.66395229482
$ws.Cells.Item(206, 8).Value = 46000
$ws.Cells.Item(206, 9).Value = 12475.08199
$ws.Cells.Item(206, 10).Value = 47.2727272727273
$ws.Cells.Item(206, 11).Value = 60
$ws.Cells.Item(206, 12).Value = 380
$ws.Cells.Item(206, 13).Value = 2760
$ws.Cells.Item(206, 14).Value = 8250.5
$ws.Cells.Item(206, 15).Value = 1816987.417
$ws.Cells.Item(206, 16).Value = 5524893.755
$ws.Cells.Item(206, 17).Value = 'Palmerston North City'
$ws.Cells.Item(206, 18).Value = 'Manawatū'
$ws.Cells.Item(206, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(206, 20).Value = 'Mana_11a'
$ws.Cells.Item(206, 21).Value = '% exceedances over 540/100 mL'
# Row 207
$ws.Cells.Item(207, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(207, 2).Value = 'E coli (Median)'
$ws.Cells.Item(207, 3).Value = 'E'
$ws.Cells.Item(207, 4).Value = '2019 - 2023'
$ws.Cells.Item(207, 5).Value = 'RepSite'
$ws.Cells.Item(207, 6).Value = 490
$ws.Cells.Item(207, 7).Value = 2624.66395229482
$ws.Cells.Item(207, 8).Value = 46000
$ws.Cells.Item(207, 9).Value = 12475.08199
$ws.Cells.Item(207, 10).Value = 47.2727272727273
$ws.Cells.Item(207, 11).Value = 60
$ws.Cells.Item(207, 12).Value = 380
$ws.Cells.Item(207, 13).Value = 2760
$ws.Cells.Item(207, 14).Value = 8250.5
$ws.Cells.Item(207, 15).Value = 1816987.417
$ws.Cells.Item(207, 16).Value = 5524893.755
$ws.Cells.Item(207, 17).Value = 'Palmerston North City'
$ws.Cells.Item(207, 18).Value = 'Manawatū'
$ws.Cells.Item(207, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(207, 20).Value = 'Mana_11a'
$ws.Cells.Item(207, 21).Value = 'E. coli/100 mL'
# Row 208
$ws.Cells.Item(208, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(208, 2).Value = 'E coli (95th Percentile)'
$ws.Cells.Item(208, 3).Value = 'E'
$ws.Cells.Item(208, 4).Value = '2019 - 2023'
$ws.Cells.Item(208, 5).Value = 'RepSite'
$ws.Cells.Item(208, 6).Value = 490
$ws.Cells.Item(208, 7).Value = 2624.66395229482
$ws.Cells.Item(208, 8).Value = 46000
$ws.Cells.Item(208, 9).Value = 12475.08199
$ws.Cells.Item(208, 10).Value = 47.2727272727273
$ws.Cells.Item(208, 11).Value = 60
$ws.Cells.Item(208, 12).Value = 380
$ws.Cells.Item(208, 13).Value = 2760
$ws.Cells.Item(208, 14).Value = 8250.5
$ws.Cells.Item(208, 15).Value = 1816987.417
$ws.Cells.Item(208, 16).Value = 5524893.755
$ws.Cells.Item(208, 17).Value = 'Palmerston North City'
$ws.Cells.Item(208, 18).Value = 'Manawatū'
$ws.Cells.Item(208, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(208, 20).Value = 'Mana_11a'
$ws.Cells.Item(208, 21).Value = 'E. coli/100 mL'
# Row 209
$ws.Cells.Item(209, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(209, 2).Value = 'Ammoniacal-N (95th Percentile)'
$ws.Cells.Item(209, 3).Value = 'B'
$ws.Cells.Item(209, 4).Value = '2019 - 2023'
$ws.Cells.Item(209, 5).Value = 'RepSite'
$ws.Cells.Item(209, 6).Value = 0.06510000000000001
$ws.Cells.Item(209, 7).Value = 0.103637813504754
$ws.Cells.Item(209, 8).Value = 0.823303787577613
$ws.Cells.Item(209, 9).Value = 0.23162
$ws.Cells.Item(209, 12).Value = 0.07543
$ws.Cells.Item(209, 13).Value = 0.14335
$ws.Cells.Item(209, 14).Value = 0.20668
$ws.Cells.Item(209, 15).Value = 1816987.417
$ws.Cells.Item(209, 16).Value = 5524893.755
$ws.Cells.Item(209, 17).Value = 'Palmerston North City'
$ws.Cells.Item(209, 18).Value = 'Manawatū'
$ws.Cells.Item(209, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(209, 20).Value = 'Mana_11a'
$ws.Cells.Item(209, 21).Value = 'mg NH4-N/L'
# Row 210
$ws.Cells.Item(210, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(210, 2).Value = 'Ammoniacal-N (Median)'
$ws.Cells.Item(210, 3).Value = 'B'
$ws.Cells.Item(210, 4).Value = '2019 - 2023'
$ws.Cells.Item(210, 5).Value = 'RepSite'
$ws.Cells.Item(210, 6).Value = 0.06510000000000001
$ws.Cells.Item(210, 7).Value = 0.103637813504754
$ws.Cells.Item(210, 8).Value = 0.823303787577613
$ws.Cells.Item(210, 9).Value = 0.23162
$ws.Cells.Item(210, 12).Value = 0.07543
$ws.Cells.Item(210, 13).Value = 0.14335
$ws.Cells.Item(210, 14).Value = 0.20668
$ws.Cells.Item(210, 15).Value = 1816987.417
$ws.Cells.Item(210, 16).Value = 5524893.755
$ws.Cells.Item(210, 17).Value = 'Palmerston North City'
$ws.Cells.Item(210, 18).Value = 'Manawatū'
$ws.Cells.Item(210, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(210, 20).Value = 'Mana_11a'
$ws.Cells.Item(210, 21).Value = 'mg NH4-N/L'
# Row 211
$ws.Cells.Item(211, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(211, 2).Value = 'Nitrate-N (95th Percentile)'
$ws.Cells.Item(211, 3).Value = 'A'
$ws.Cells.Item(211, 4).Value = '2019 - 2023'
$ws.Cells.Item(211, 5).Value = 'RepSite'
$ws.Cells.Item(211, 6).Value = 0.412
$ws.Cells.Item(211, 7).Value = 0.460509090909091
$ws.Cells.Item(211, 8).Value = 1.04
$ws.Cells.Item(211, 9).Value = 0.959
$ws.Cells.Item(211, 12).Value = 0.24
$ws.Cells.Item(211, 13).Value = 0.77615
$ws.Cells.Item(211, 14).Value = 0.9011
$ws.Cells.Item(211, 15).Value = 1816987.417
$ws.Cells.Item(211, 16).Value = 5524893.755
$ws.Cells.Item(211, 17).Value = 'Palmerston North City'
$ws.Cells.Item(211, 18).Value = 'Manawatū'
$ws.Cells.Item(211, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(211, 20).Value = 'Mana_11a'
$ws.Cells.Item(211, 21).Value = 'mg NO3-N/L'
# Row 212
$ws.Cells.Item(212, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(212, 2).Value = 'Nitrate-N (Median)'
$ws.Cells.Item(212, 3).Value = 'A'
$ws.Cells.Item(212, 4).Value = '2019 - 2023'
$ws.Cells.Item(212, 5).Value = 'RepSite'
$ws.Cells.Item(212, 6).Value = 0.412
$ws.Cells.Item(212, 7).Value = 0.460509090909091
$ws.Cells.Item(212, 8).Value = 1.04
$ws.Cells.Item(212, 9).Value = 0.959
$ws.Cells.Item(212, 12).Value = 0.24
$ws.Cells.Item(212, 13).Value = 0.77615
$ws.Cells.Item(212, 14).Value = 0.9011
$ws.Cells.Item(212, 15).Value = 1816987.417
$ws.Cells.Item(212, 16).Value = 5524893.755
$ws.Cells.Item(212, 17).Value = 'Palmerston North City'
$ws.Cells.Item(212, 18).Value = 'Manawatū'
$ws.Cells.Item(212, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(212, 20).Value = 'Mana_11a'
$ws.Cells.Item(212, 21).Value = 'mg NO3-N/L'
# Row 213
$ws.Cells.Item(213, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(213, 2).Value = 'Soluble Inorganic Nitrogen (95th Percentile)'
$ws.Cells.Item(213, 4).Value = '2019 - 2023'
$ws.Cells.Item(213, 5).Value = 'RepSite'
$ws.Cells.Item(213, 6).Value = 0.556
$ws.Cells.Item(213, 7).Value = 0.595181818181818
$ws.Cells.Item(213, 8).Value = 1.082
$ws.Cells.Item(213, 9).Value = 1.04925
$ws.Cells.Item(213, 12).Value = 0.394
$ws.Cells.Item(213, 13).Value = 0.8596
$ws.Cells.Item(213, 14).Value = 0.993
$ws.Cells.Item(213, 15).Value = 1816987.417
$ws.Cells.Item(213, 16).Value = 5524893.755
$ws.Cells.Item(213, 17).Value = 'Palmerston North City'
$ws.Cells.Item(213, 18).Value = 'Manawatū'
$ws.Cells.Item(213, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(213, 20).Value = 'Mana_11a'
$ws.Cells.Item(213, 21).Value = 'g/m3'
# Row 214
$ws.Cells.Item(214, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(214, 2).Value = 'Soluble Inorganic Nitrogen (Median)'
$ws.Cells.Item(214, 4).Value = '2019 - 2023'
$ws.Cells.Item(214, 5).Value = 'RepSite'
$ws.Cells.Item(214, 6).Value = 0.556
$ws.Cells.Item(214, 7).Value = 0.595181818181818
$ws.Cells.Item(214, 8).Value = 1.082
$ws.Cells.Item(214, 9).Value = 1.04925
$ws.Cells.Item(214, 12).Value = 0.394
$ws.Cells.Item(214, 13).Value = 0.8596
$ws.Cells.Item(214, 14).Value = 0.993
$ws.Cells.Item(214, 15).Value = 1816987.417
$ws.Cells.Item(214, 16).Value = 5524893.755
$ws.Cells.Item(214, 17).Value = 'Palmerston North City'
$ws.Cells.Item(214, 18).Value = 'Manawatū'
$ws.Cells.Item(214, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(214, 20).Value = 'Mana_11a'
$ws.Cells.Item(214, 21).Value = 'g/m3'
# Row 215
$ws.Cells.Item(215, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(215, 2).Value = 'Total Nitrogen (95th Percentile)'
$ws.Cells.Item(215, 4).Value = '2019 - 2023'
$ws.Cells.Item(215, 5).Value = 'RepSite'
$ws.Cells.Item(215, 6).Value = 0.85
$ws.Cells.Item(215, 7).Value = 0.887818181818182
$ws.Cells.Item(215, 8).Value = 1.65
$ws.Cells.Item(215, 9).Value = 1.4
$ws.Cells.Item(215, 12).Value = 0.74
$ws.Cells.Item(215, 13).Value = 1.1945
$ws.Cells.Item(215, 14).Value = 1.331
$ws.Cells.Item(215, 15).Value = 1816987.417
$ws.Cells.Item(215, 16).Value = 5524893.755
$ws.Cells.Item(215, 17).Value = 'Palmerston North City'
$ws.Cells.Item(215, 18).Value = 'Manawatū'
$ws.Cells.Item(215, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(215, 20).Value = 'Mana_11a'
$ws.Cells.Item(215, 21).Value = 'g/m3'
# Row 216
$ws.Cells.Item(216, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(216, 2).Value = 'Total Nitrogen (Median)'
$ws.Cells.Item(216, 4).Value = '2019 - 2023'
$ws.Cells.Item(216, 5).Value = 'RepSite'
$ws.Cells.Item(216, 6).Value = 0.85
$ws.Cells.Item(216, 7).Value = 0.887818181818182
$ws.Cells.Item(216, 8).Value = 1.65
$ws.Cells.Item(216, 9).Value = 1.4
$ws.Cells.Item(216, 12).Value = 0.74
$ws.Cells.Item(216, 13).Value = 1.1945
$ws.Cells.Item(216, 14).Value = 1.331
$ws.Cells.Item(216, 15).Value = 1816987.417
$ws.Cells.Item(216, 16).Value = 5524893.755
$ws.Cells.Item(216, 17).Value = 'Palmerston North City'
$ws.Cells.Item(216, 18).Value = 'Manawatū'
$ws.Cells.Item(216, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(216, 20).Value = 'Mana_11a'
$ws.Cells.Item(216, 21).Value = 'g/m3'
# Row 217
$ws.Cells.Item(217, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(217, 2).Value = 'Total Phosphorus (95th Percentile)'
$ws.Cells.Item(217, 4).Value = '2019 - 2023'
$ws.Cells.Item(217, 5).Value = 'RepSite'
$ws.Cells.Item(217, 6).Value = 0.058
$ws.Cells.Item(217, 7).Value = 0.112745454545455
$ws.Cells.Item(217, 8).Value = 0.755
$ws.Cells.Item(217, 9).Value = 0.39425
$ws.Cells.Item(217, 12).Value = 0.052
$ws.Cells.Item(217, 13).Value = 0.2013
$ws.Cells.Item(217, 14).Value = 0.2657
$ws.Cells.Item(217, 15).Value = 1816987.417
$ws.Cells.Item(217, 16).Value = 5524893.755
$ws.Cells.Item(217, 17).Value = 'Palmerston North City'
$ws.Cells.Item(217, 18).Value = 'Manawatū'
$ws.Cells.Item(217, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(217, 20).Value = 'Mana_11a'
$ws.Cells.Item(217, 21).Value = 'g/m3'
# Row 218
$ws.Cells.Item(218, 1).Value = 'Manawatu at us Fonterra Longburn'
$ws.Cells.Item(218, 2).Value = 'Total Phosphorus (Median)'
$ws.Cells.Item(218, 4).Value = '2019 - 2023'
$ws.Cells.Item(218, 5).Value = 'RepSite'
$ws.Cells.Item(218, 6).Value = 0.058
$ws.Cells.Item(218, 7).Value = 0.112745454545455
$ws.Cells.Item(218, 8).Value = 0.755
$ws.Cells.Item(218, 9).Value = 0.39425
$ws.Cells.Item(218, 12).Value = 0.052
$ws.Cells.Item(218, 13).Value = 0.2013
$ws.Cells.Item(218, 14).Value = 0.2657
$ws.Cells.Item(218, 15).Value = 1816987.417
$ws.Cells.Item(218, 16).Value = 5524893.755
$ws.Cells.Item(218, 17).Value = 'Palmerston North City'
$ws.Cells.Item(218, 18).Value = 'Manawatū'
$ws.Cells.Item(218, 19).Value = 'Lower Manawatu'
$ws.Cells.Item(218, 20).Value = 'Mana_11a'
$ws.Cells.Item(218, 21).Value = 'g/m3'

Write-Host "Edit applied: 150 scalar updates + 17 new rows (202-218)."
